$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Row 10 (No.) / Row 11 (Marking) / Row 12 (Total) - update counts / labels
# ---------------------------------------------------------------------------

# Give the row-label cells (A10/A11/A12) the same bold "mtitleStyle" (s=4)
# already used by A9:E9 / row 15 headers, by copying format from A9.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)

# Row 10: Right / Wrong / Not-attempted / Max counts
$ws.Range("B10").Value = 10
$ws.Range("D10").Value = 18
$ws.Range("E10").Value = 28

# Row 11: Marking scheme - fix "-1" which was stored as text so it behaves
# as a real number (the float-input bug referenced by the commit message).
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12: Totals
$ws.Range("B12").Value = 40
$ws.Range("E12").Value = "40/112"

# ---------------------------------------------------------------------------
# Remove the third (G:H) "Student Ans / Correct Ans" column pair entirely.
# ---------------------------------------------------------------------------
$ws.Columns("G:H").Delete()

# Remove the second (D:E) "Student Ans / Correct Ans" pair for every question
# row except the first three (16-18), which are kept (and updated below).
$ws.Range("D19:E40").Clear()

# ---------------------------------------------------------------------------
# Populate the "Student Ans" cells (column A, and D for rows 16-18) that now
# match the "Correct Ans" column, using the same green "correctStyle" (s=5)
# already applied to B10:B12.
# ---------------------------------------------------------------------------
$ws.Range("B10").Copy()

$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = "Option A"

$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = "Option C"

$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = "Option B"
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "Option D"

$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A21").Value = "Option C"

$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("A25").Value = "Option A"

$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("A30").Value = "Option B"

$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("A33").Value = "Option D"

$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("A36").Value = "Option A"

$ws.Range("A39").PasteSpecial(-4122)
$ws.Range("A39").Value = "Option D"
